# Update the LR-pair results table on Sheet1 with the full 3x4 sender/target
# cluster combination (ECs/FAPs/MuSCs senders x ECs/FAPs/MuSCs/Resolving-Mac
# targets) for Efna5 -> Epha3, replacing the previous partial 3x3 table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("ECs", "Efna5", "Epha3", "ECs", 3, 1, 0.2708946666666667, 0.812684, 0.1616296696421007, 0.1616296696421007, 1, 0.3333333333333333, 0.009849666666666666, 0.029549, 0.0002013876315934659, 0.0002013876315934659, 0.002668222168444444, 0.024013999516, 0.00003255021636445697, 0.00003255021636445697),
    @("ECs", "Efna5", "Epha3", "FAPs", 3, 1, 0.2708946666666667, 0.812684, 0.1616296696421007, 0.1616296696421007, 3, 1, 46.25093466666667, 138.752804, 0.9456529349389956, 0.9456529349389956, 12.52913152954844, 112.762183765936, 0.1528455714702728, 0.1528455714702728),
    @("ECs", "Efna5", "Epha3", "MuSCs", 3, 1, 0.2708946666666667, 0.812684, 0.1616296696421007, 0.1616296696421007, 3, 1, 2.607896333333333, 7.823689, 0.05332140505715427, 0.05332140505715428, 0.7064652079195555, 6.358186871276, 0.008618321084240479, 0.008618321084240481),
    @("ECs", "Efna5", "Epha3", "Resolving-Mac", 3, 1, 0.2708946666666667, 0.812684, 0.1616296696421007, 0.1616296696421007, 1, 0.3333333333333333, 0.04031433333333333, 0.120943, 0.0008242723722565416, 0.0008242723722565415, 0.01092093789022222, 0.09828844101199999, 0.0001332268712229354, 0.0001332268712229354),
    @("FAPs", "Efna5", "Epha3", "ECs", 2, 0.6666666666666666, 1.058121333333333, 3.174364, 0.6313295261673385, 0.6313295261673384, 1, 0.3333333333333333, 0.009849666666666666, 0.029549, 0.0002013876315934659, 0.0002013876315934659, 0.01042214242622222, 0.093799281836, 0.0001271419580298654, 0.0001271419580298653),
    @("FAPs", "Efna5", "Epha3", "FAPs", 2, 0.6666666666666666, 1.058121333333333, 3.174364, 0.6313295261673385, 0.6313295261673384, 3, 1, 46.25093466666667, 138.752804, 0.9456529349389956, 0.9456529349389956, 48.93910065740623, 440.451905916656, 0.5970186193337891, 0.597018619333789),
    @("FAPs", "Efna5", "Epha3", "MuSCs", 2, 0.6666666666666666, 1.058121333333333, 3.174364, 0.6313295261673385, 0.6313295261673384, 3, 1, 2.607896333333333, 7.823689, 0.05332140505715427, 0.05332140505715428, 2.759470745421778, 24.835236708796, 0.03366337738930993, 0.03366337738930993),
    @("FAPs", "Efna5", "Epha3", "Resolving-Mac", 2, 0.6666666666666666, 1.058121333333333, 3.174364, 0.6313295261673385, 0.6313295261673384, 1, 0.3333333333333333, 0.04031433333333333, 0.120943, 0.0008242723722565416, 0.0008242723722565415, 0.04265745613911111, 0.383917105252, 0.0005203874862095505, 0.0005203874862095504),
    @("MuSCs", "Efna5", "Epha3", "ECs", 2, 0.6666666666666666, 0.3470046666666667, 1.041014, 0.2070408041905609, 0.2070408041905609, 1, 0.3333333333333333, 0.009849666666666666, 0.029549, 0.0002013876315934659, 0.0002013876315934659, 0.003417880298444444, 0.030760922686, 0.00004169545719914359, 0.00004169545719914359),
    @("MuSCs", "Efna5", "Epha3", "FAPs", 2, 0.6666666666666666, 0.3470046666666667, 1.041014, 0.2070408041905609, 0.2070408041905609, 3, 1, 46.25093466666667, 138.752804, 0.9456529349389956, 0.9456529349389956, 16.04929016702845, 144.443611503256, 0.1957887441349338, 0.1957887441349338),
    @("MuSCs", "Efna5", "Epha3", "MuSCs", 2, 0.6666666666666666, 0.3470046666666667, 1.041014, 0.2070408041905609, 0.2070408041905609, 3, 1, 2.607896333333333, 7.823689, 0.05332140505715427, 0.05332140505715428, 0.9049521978495555, 8.144569780646, 0.01103970658360386, 0.01103970658360386),
    @("MuSCs", "Efna5", "Epha3", "Resolving-Mac", 2, 0.6666666666666666, 0.3470046666666667, 1.041014, 0.2070408041905609, 0.2070408041905609, 1, 0.3333333333333333, 0.04031433333333333, 0.120943, 0.0008242723722565416, 0.0008242723722565415, 0.01398926180022222, 0.125903356202, 0.0001706580148240557, 0.0001706580148240557),
)

$r = 2
foreach ($row in $rows) {
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $row[$c]
    }
    $r = $r + 1
}

